# Updates loading_percent results (rows 2-25, columns B,D,E,F,G,H,I,J,K,M)
# for the "case with 380 kV" re-run — recomputed loading percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: "row,col,newValue"
$data = @(
"2,2,5.105260970161016",
"2,4,7.734427766164324",
"2,5,12.87115848005571",
"2,6,39.06502231613246",
"2,7,45.34960192973936",
"2,8,18.3360737790213",
"2,9,22.97514253333708",
"2,10,10.16167643963673",
"2,11,13.37848043877065",
"2,13,16.7779707221353",
"3,2,4.994088738418792",
"3,4,7.724163839576789",
"3,5,12.8825408235109",
"3,6,39.08339835267538",
"3,7,45.31163704500311",
"3,8,18.37968492086603",
"3,9,23.08201075881392",
"3,10,10.18495694282955",
"3,11,13.09973595293932",
"3,13,16.68250135617813",
"4,2,4.925341264764672",
"4,4,7.718738657847159",
"4,5,12.89144942798554",
"4,6,39.10458166608962",
"4,7,45.30242945218601",
"4,8,18.41001702376602",
"4,9,23.15136510887202",
"4,10,10.20039090674113",
"4,11,12.9282999792479",
"4,13,16.62669500945713",
"5,2,4.897243971181653",
"5,4,7.716749998288769",
"5,5,12.89556282185804",
"5,6,39.11569928346631",
"5,7,45.30222053793109",
"5,8,18.42326958473431",
"5,9,23.18056821681583",
"5,10,10.20696722296785",
"5,11,12.85846688762016",
"5,13,16.60467866676762",
"6,2,4.892574683162796",
"6,4,7.716433245879796",
"6,5,12.89627503099881",
"6,6,39.11769530660366",
"6,7,45.30239969787436",
"6,8,18.42552397911055",
"6,9,23.18547421318016",
"6,10,10.20807654853772",
"6,11,12.84687585186005",
"6,13,16.60106717666826",
"7,2,4.92496261302518",
"7,4,7.71871093652729",
"7,5,12.89150294639367",
"7,6,39.10472154624557",
"7,7,45.30241229479907",
"7,8,18.41019214343247",
"7,9,23.15175514224249",
"7,10,10.20047843542701",
"7,11,12.92735793627412",
"7,13,16.62639513002516",
"8,2,5.067053703045936",
"8,4,7.730707773057449",
"8,5,12.87468487433356",
"8,6,39.06930137497289",
"8,7,45.33358369868623",
"8,8,18.35037203086102",
"8,9,23.01121571414357",
"8,10,10.16946719254075",
"8,11,13.28248803776129",
"8,13,16.74448033473434",
"9,2,5.340130918749437",
"9,4,7.761116152508832",
"9,5,12.85692086686935",
"9,6,39.078533025732",
"9,7,45.50660194888193",
"9,8,18.26135008574642",
"9,9,22.76522252690975",
"9,10,10.11768441090147",
"9,11,13.97236576097373",
"9,13,16.99752107075498",
"10,2,5.556915723721876",
"10,4,7.787547394938009",
"10,5,12.85311952728363",
"10,6,39.13339074053011",
"10,7,45.70170507620278",
"10,8,18.21329364165721",
"10,9,22.60247397107912",
"10,10,10.08512782142449",
"10,11,14.46978060538678",
"10,13,17.19534078114607",
"11,2,5.670304770050865",
"11,4,7.800435282214453",
"11,5,12.85339097420499",
"11,6,39.16877792949015",
"11,7,45.80510529119837",
"11,8,18.19521845484358",
"11,9,22.53232711612117",
"11,10,10.0715051703445",
"11,11,14.6929493969851",
"11,13,17.28765265276818",
"12,2,5.712503237226041",
"12,4,7.805437593479793",
"12,5,12.85378057908316",
"12,6,39.18367488556891",
"12,7,45.84635064031978",
"12,8,18.18891966758415",
"12,9,22.50632265381231",
"12,10,10.06651712219845",
"12,11,14.77692833383226",
"12,13,17.32291923259733",
"13,2,5.703448226492369",
"13,4,7.80435487055292",
"13,5,12.85368392993177",
"13,6,39.18040007081311",
"13,7,45.83737506155005",
"13,8,18.19025192331837",
"13,9,22.51189834228795",
"13,10,10.0675838054696",
"13,11,14.75886682051335",
"13,13,17.31531053881209",
"14,2,5.67379139074601",
"14,4,7.800844394986576",
"14,5,12.85341728454491",
"14,6,39.16997354279139",
"14,7,45.80845679090667",
"14,8,18.19468930134195",
"14,9,22.53017652028274",
"14,10,10.07109138364859",
"14,11,14.69986948005334",
"14,13,17.29054797982982",
"15,2,5.655528861363266",
"15,4,7.798709938056188",
"15,5,12.85329127989688",
"15,6,39.16378176285572",
"15,7,45.79101514462272",
"15,8,18.19747845643606",
"15,9,22.54144516810053",
"15,10,10.07326208160111",
"15,11,14.66366046914441",
"15,13,17.27541986043734",
"16,2,5.549403041606907",
"16,4,7.786722334639825",
"16,5,12.85314197388928",
"16,6,39.13128775925666",
"16,7,45.69524104894511",
"16,8,18.21455121020487",
"16,9,22.60713643120046",
"16,10,10.086041971581",
"16,11,14.45512621874892",
"16,13,17.18935268817072",
"17,2,5.483001937386691",
"17,4,7.779588046322894",
"17,5,12.85356221482476",
"17,6,39.11402364103056",
"17,7,45.64022767169257",
"17,8,18.22599555969322",
"17,9,22.64843140295151",
"17,10,10.09418602067636",
"17,11,14.32634027182437",
"17,13,17.13713130484985",
"18,2,5.449956590363348",
"18,4,7.775566078488792",
"18,5,12.85399230167397",
"18,6,39.10507600894829",
"18,7,45.60996525484333",
"18,8,18.23293431689193",
"18,9,22.67254918842936",
"18,10,10.09898204822986",
"18,11,14.2519781442584",
"18,13,17.1073150362058",
"19,2,5.440066166913179",
"19,4,7.774218373690411",
"19,5,12.85417030030481",
"19,6,39.10221527608492",
"19,7,45.59995634549211",
"19,8,18.23534480310467",
"19,9,22.68077793486652",
"19,10,10.10062510643347",
"19,11,14.22675359209361",
"19,13,17.09725828123635",
"20,2,5.490119222427582",
"20,4,7.780339086495282",
"20,5,12.85349798934061",
"20,6,39.11575979913475",
"20,7,45.64594123991954",
"20,8,18.22474040529904",
"20,9,22.64399760581128",
"20,10,10.09330750482655",
"20,11,14.3400801684477",
"20,13,17.14266775348599",
"21,2,5.682522540564344",
"21,4,7.801872216093579",
"21,5,12.85348782835255",
"21,6,39.17299548549498",
"21,7,45.81689421445576",
"21,8,18.19337110937498",
"21,9,22.52479261903462",
"21,10,10.07005649577183",
"21,11,14.71721345761936",
"21,13,17.29781312126695",
"22,2,5.803951243760771",
"22,4,7.816655023669392",
"22,5,12.85515246312368",
"22,6,39.21912288450719",
"22,7,45.94079363848707",
"22,8,18.17605192972983",
"22,9,22.45014119200209",
"22,10,10.05585464163872",
"22,11,14.9605663196908",
"22,13,17.40100554571169",
"23,2,5.739543577881943",
"23,4,7.808701018294344",
"23,5,12.85411141739676",
"23,6,39.19370745958923",
"23,7,45.87355872192212",
"23,8,18.18500385464338",
"23,9,22.4896862986476",
"23,10,10.06334355357713",
"23,11,14.83099646502895",
"23,13,17.34577349806842",
"24,2,5.486903013341864",
"24,4,7.77999929303966",
"24,5,12.85352643851334",
"24,6,39.11497183609418",
"24,7,45.64335388092519",
"24,8,18.22530674160609",
"24,9,22.64600095183512",
"24,10,10.09370432701467",
"24,11,14.33386935576035",
"24,13,17.14016407999211",
"25,2,5.267081067559073",
"25,4,7.752163538688479",
"25,5,12.86009977742192",
"25,6,39.06759476070308",
"25,7,45.44782674795888",
"25,8,18.28239346663744",
"25,9,22.82860774931985",
"25,10,10.13072802608274",
"25,11,13.78701431641618",
"25,13,16.92688639899168"
)

foreach ($line in $data) {
    $parts = $line -split ','
    $r = [int]$parts[0]
    $c = [int]$parts[1]
    $v = [double]$parts[2]
    $ws.Cells.Item($r, $c).Value = $v
}
